$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $newText, 2)
}

Replace-ExactText "-0.021" "0.014"
Replace-ExactText "0.259" "0.335"
Replace-ExactText "0.002" "0.017"
Replace-ExactText "0.517" "0.36"
Replace-ExactText "0.012" "-0.010"
Replace-ExactText "0.776" "0.663"
